# Remake the hero skill and relative card for defender (铁盾 / Iron Shield)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HeroSkill")

# Row 5 corresponds to the "铁盾" (Iron Shield) hero skill / card entry.
# Its description, type, and linked CardId are reworked:
#   Des:    可以给己方单位装备一面铁盾  ->  使我方王塔获得一些物甲
#   Type:   2                         ->  3
#   CardId: 52100001                  ->  53100008
$ws.Range("C5").Value = "使我方王塔获得一些物甲"
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 53100008

# Update the selected/active cell to match the edited row.
$ws.Range("E5").Select()
